$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column C (Algorithm result column) for specific rows
# Full-precision decimal literals are used so the stored IEEE-754 double
# exactly matches the source data (some of these are floating point
# artifacts, e.g. -11.86359999999999 rather than the "clean" -11.8636).
$ws.Range("C10").Value = -13.7439
$ws.Range("C12").Value = -10.4219
$ws.Range("C18").Value = -11.86359999999999
$ws.Range("C37").Value = -12.9227
$ws.Range("C55").Value = -13.43009999999999
$ws.Range("C68").Value = -11.6173
$ws.Range("C77").Value = -13.3782
$ws.Range("C78").Value = -13.1528
